$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q4").Value = 3.4
$ws.Range("R4").Value = 1.33
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.25
$ws.Range("K5").Value = 2.2
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 9.5
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 3.4
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.75
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("X5").Value = 7
$ws.Range("Z5").Value = 12
$ws.Range("AC5").Value = 9.5
$ws.Range("AE5").Value = 19
$ws.Range("AG5").Value = 13
$ws.Range("AL5").Value = 51
$ws.Range("AN5").Value = 3.5
$ws.Range("AO5").Value = 8.5
$ws.Range("AT5").Value = 2.75
$ws.Range("AU5").Value = 9
$ws.Range("AW5").Value = 7
$ws.Range("AY5").Value = 41
$ws.Range("G6").Value = 2.3
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 3.4
$ws.Range("J6").Value = 3.25
$ws.Range("L6").Value = 4.33
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 2.25
$ws.Range("S6").Value = 1.62
$ws.Range("T6").Value = 2.2
$ws.Range("W6").Value = 5.5
$ws.Range("X6").Value = 9.5
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 26
$ws.Range("AG6").Value = 7
$ws.Range("AH6").Value = 15
$ws.Range("AJ6").Value = 41
$ws.Range("AN6").Value = 4
$ws.Range("AO6").Value = 15
$ws.Range("AT6").Value = 2.2
$ws.Range("AW6").Value = 5
$ws.Range("AX6").Value = 21
$ws.Range("AY6").Value = 41
$ws.Range("AZ6").Value = 81
$ws.Range("S7").Value = 1.67
$ws.Range("T7").Value = 2.1
$ws.Range("J12").Value = 2.27
$ws.Range("AT12").Value = 3
$ws.Range("AV12").Value = 55
$ws.Range("G15").Value = 1.7
$ws.Range("H15").Value = 3.3
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 2.4
$ws.Range("Q15").Value = 2.35
$ws.Range("R15").Value = 1.57
$ws.Range("S15").Value = 1.5
$ws.Range("T15").Value = 2.5
$ws.Range("AD15").Value = 7
$ws.Range("AG15").Value = 11
$ws.Range("AN15").Value = 3.5
$ws.Range("AO15").Value = 9.5
$ws.Range("AT15").Value = 2.5
$ws.Range("AZ15").Value = 126
$ws.Range("G25").Value = 2
$ws.Range("I25").Value = 3.8
$ws.Range("J25").Value = 2.63
$ws.Range("L25").Value = 4
$ws.Range("M25").Value = 1.05
$ws.Range("N25").Value = 11
$ws.Range("Q25").Value = 1.98
$ws.Range("R25").Value = 1.88
$ws.Range("U25").Value = 1.73
$ws.Range("V25").Value = 2
$ws.Range("X25").Value = 9.5
$ws.Range("AA25").Value = 17
$ws.Range("AC25").Value = 10
$ws.Range("AE25").Value = 13
$ws.Range("AF25").Value = 41
$ws.Range("AG25").Value = 11
$ws.Range("AH25").Value = 19
$ws.Range("AK25").Value = 29
$ws.Range("AL25").Value = 34
$ws.Range("AM25").Value = 201
$ws.Range("AQ25").Value = 41
$ws.Range("G26").Value = 1.4
$ws.Range("H26").Value = 5.5
$ws.Range("Q26").Value = 1.33
$ws.Range("R26").Value = 3.4
$ws.Range("S26").Value = 1.2
$ws.Range("T26").Value = 4.33
$ws.Range("AD26").Value = 13
$ws.Range("AG26").Value = 26
$ws.Range("AI26").Value = 19
$ws.Range("AZ26").Value = 81
$ws.Range("BD26").Value = 151
$ws.Range("G28").Value = 2.37
$ws.Range("H28").Value = 2.9
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 2.02
$ws.Range("L28").Value = 3.6
$ws.Range("M28").Value = 1.09
$ws.Range("N28").Value = 6.3
$ws.Range("O28").Value = 1.37
$ws.Range("P28").Value = 2.85
$ws.Range("Q28").Value = 2.15
$ws.Range("R28").Value = 1.65
$ws.Range("S28").Value = 1.44
$ws.Range("T28").Value = 2.62
$ws.Range("W28").Value = 7.3
$ws.Range("X28").Value = 11.5
$ws.Range("Y28").Value = 9.25
$ws.Range("Z28").Value = 26
$ws.Range("AA28").Value = 21
$ws.Range("AB28").Value = 32
$ws.Range("AC28").Value = 6.3
$ws.Range("AD28").Value = 5.7
$ws.Range("AE28").Value = 13.5
$ws.Range("AF28").Value = 65
$ws.Range("AG28").Value = 8.5
$ws.Range("AI28").Value = 10.75
$ws.Range("AL28").Value = 35
$ws.Range("AN28").Value = 4.35
$ws.Range("AO28").Value = 13
$ws.Range("AP28").Value = 20
$ws.Range("AQ28").Value = 55
$ws.Range("AR28").Value = 80
$ws.Range("AT28").Value = 2.62
$ws.Range("AU28").Value = 6.7
$ws.Range("AY28").Value = 23
$ws.Range("BA28").Value = 110
$ws.Range("G29").Value = 2.8
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 2.18
$ws.Range("J29").Value = 3.35
$ws.Range("K29").Value = 2.22
$ws.Range("L29").Value = 2.77
$ws.Range("N29").Value = 8
$ws.Range("O29").Value = 1.25
$ws.Range("S29").Value = 1.35
$ws.Range("T29").Value = 2.92
$ws.Range("W29").Value = 10.25
$ws.Range("X29").Value = 15.5
$ws.Range("Y29").Value = 10.5
$ws.Range("Z29").Value = 32
$ws.Range("AA29").Value = 23
$ws.Range("AB29").Value = 29
$ws.Range("AC29").Value = 8
$ws.Range("AE29").Value = 13.5
$ws.Range("AH29").Value = 11.5
$ws.Range("AI29").Value = 9
$ws.Range("AJ29").Value = 21
$ws.Range("AK29").Value = 16.5
$ws.Range("AN29").Value = 4.9
$ws.Range("AO29").Value = 15
$ws.Range("AP29").Value = 22
$ws.Range("AQ29").Value = 65
$ws.Range("AR29").Value = 100
$ws.Range("AT29").Value = 2.92
$ws.Range("AU29").Value = 7
$ws.Range("AW29").Value = 4.25
$ws.Range("AX29").Value = 11.25
$ws.Range("AY29").Value = 18.5
$ws.Range("BA29").Value = 70
$ws.Range("BB29").Value = 200
$ws.Range("G30").Value = 1.75
$ws.Range("H30").Value = 3.85
$ws.Range("I30").Value = 3.75
$ws.Range("J30").Value = 2.22
$ws.Range("K30").Value = 2.42
$ws.Range("L30").Value = 4
$ws.Range("P30").Value = 4.6
$ws.Range("R30").Value = 2.42
$ws.Range("S30").Value = 1.26
$ws.Range("T30").Value = 3.45
$ws.Range("U30").Value = 1.5
$ws.Range("V30").Value = 2.4
$ws.Range("W30").Value = 10.75
$ws.Range("X30").Value = 10.75
$ws.Range("Z30").Value = 15.5
$ws.Range("AA30").Value = 12
$ws.Range("AB30").Value = 18
$ws.Range("AD30").Value = 8.25
$ws.Range("AE30").Value = 12.5
$ws.Range("AF30").Value = 40
$ws.Range("AG30").Value = 16
$ws.Range("AH30").Value = 25
$ws.Range("AI30").Value = 13
$ws.Range("AJ30").Value = 55
$ws.Range("AK30").Value = 29
$ws.Range("AL30").Value = 29
$ws.Range("AN30").Value = 4.05
$ws.Range("AO30").Value = 8.5
$ws.Range("AP30").Value = 13.5
$ws.Range("AQ30").Value = 25
$ws.Range("AR30").Value = 40
$ws.Range("AT30").Value = 3.45
$ws.Range("AU30").Value = 6.6
$ws.Range("AW30").Value = 6.2
$ws.Range("AX30").Value = 19.5
$ws.Range("AY30").Value = 21
$ws.Range("AZ30").Value = 90
$ws.Range("BA30").Value = 100
$ws.Range("BB30").Value = 200
$ws.Range("H31").Value = 3.6
$ws.Range("I31").Value = 4.35
$ws.Range("J31").Value = 2.27
$ws.Range("L31").Value = 4.75
$ws.Range("P31").Value = 3.3
$ws.Range("Q31").Value = 1.87
$ws.Range("R31").Value = 1.87
$ws.Range("U31").Value = 1.83
$ws.Range("V31").Value = 1.88
$ws.Range("X31").Value = 7.8
$ws.Range("AE31").Value = 16
$ws.Range("AF31").Value = 75
$ws.Range("AG31").Value = 12.5
$ws.Range("AH31").Value = 25
$ws.Range("AJ31").Value = 75
$ws.Range("AL31").Value = 45
$ws.Range("AU31").Value = 7.6
$ws.Range("AW31").Value = 6.2
$ws.Range("BA31").Value = 175
$ws.Range("J33").Value = 3.4
$ws.Range("K33").Value = 2.25
$ws.Range("T33").Value = 2.95
$ws.Range("X33").Value = 15.5
$ws.Range("AI33").Value = 9
$ws.Range("AT33").Value = 2.95
$ws.Range("AV33").Value = 55
$ws.Range("AW33").Value = 4.2
$ws.Range("AX33").Value = 10.75
$ws.Range("AY33").Value = 18
$ws.Range("J37").Value = 2.95
$ws.Range("L37").Value = 3.65
$ws.Range("M37").Value = 1.03
$ws.Range("N37").Value = 7
$ws.Range("W37").Value = 6.9
$ws.Range("X37").Value = 11
$ws.Range("Y37").Value = 9.25
$ws.Range("AA37").Value = 21
$ws.Range("AB37").Value = 32
$ws.Range("AC37").Value = 7.7
$ws.Range("AG37").Value = 8.75
$ws.Range("AH37").Value = 16.5
$ws.Range("AI37").Value = 11
$ws.Range("AK37").Value = 29
$ws.Range("AL37").Value = 37
$ws.Range("AO37").Value = 12.5
$ws.Range("AP37").Value = 21
$ws.Range("AQ37").Value = 55
$ws.Range("AR37").Value = 90
$ws.Range("AX37").Value = 17
$ws.Range("AY37").Value = 23
$ws.Range("AZ37").Value = 80
$ws.Range("BA37").Value = 110
$ws.Range("BB37").Value = 300
